$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A44").Value = " как погода?"
$ws.Range("B44").Value = "Вопрос"
$ws.Range("C44").Value = 1
